$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 12 (shifts existing rows 12..73 down to 13..74)
$ws.Rows("12:12").Insert()

# The newly inserted row 12 is blank; start by copying the row that is now
# at row 13 (the data that used to live in row 12) into row 12, then
# overwrite the handful of cells that actually carry new data for the
# newly added weekly record.
for ($c = 1; $c -le 20; $c++) {
    $ws.Cells.Item(12, $c).Value = $ws.Cells.Item(13, $c).Value()
}

# New weekly record values for row 12
$ws.Range("D12").Value = 44607
$ws.Range("M12").Value = 60
$ws.Range("N12").Value = 7000
$ws.Range("O12").Value = 7500
$ws.Range("P12").Value = 7250
$ws.Range("S12").Value = 1812
